# Make a separate column for TIMEVAL in variables1 (Close #265)
#
# The "Variables" worksheet has a table (Table5) with columns:
#   pivot | order | variable-code | variable-type | fo_variable-label | fo_elimination | fo_note
# A new "timeval" column is inserted right after "variable-type" (so it
# becomes column E), pushing fo_variable-label / fo_elimination / fo_note
# one column to the right (F / G / H). The row that describes the "time"
# variable (row 3) used to flag itself via variable-type = "TIME"; that
# text value is replaced by a boolean TRUE in the new "timeval" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")
$lo = $ws.ListObjects.Item(1)

$lastRow = 5
$oldLastCol = 7   # A..G
$newTimevalCol = 5   # column E

# 1) Snapshot the current contents (value + a boolean flag) of the columns
#    that sit at/after the insertion point (old E, F, G == new F, G, H),
#    before we touch anything.
$snapshot = @{}
for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = $newTimevalCol; $c -le $oldLastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $key = "$r,$c"
        $snapshot[$key] = $cell.Value()
    }
}

# 2) Grow the table by one column (always appended at the end by this
#    object model), giving us a blank column H to play with while keeping
#    every other existing column untouched.
$lo.ListColumns.Add() | Out-Null

# 3) Shift the snapshotted old E/F/G values two+ one columns to the right
#    so they land in F/G/H.
for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = $oldLastCol; $c -ge $newTimevalCol; $c--) {
        $key = "$r,$c"
        $destCol = $c + 1
        $val = $snapshot[$key]
        if ($null -ne $val -and $val -ne "") {
            $ws.Cells.Item($r, $destCol).Value = $val
        } else {
            $ws.Cells.Item($r, $destCol).Value = ""
        }
    }
}

# 4) Clear out the (now vacated) new "timeval" column and populate it:
#    header in row 1, and TRUE for the "time" variable row (row 3).
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $newTimevalCol).Value = ""
}
$ws.Cells.Item(1, $newTimevalCol).Value = "timeval"
$ws.Cells.Item(3, $newTimevalCol).Value = $true

# 5) The old variable-type cell for the "time" row held the literal text
#    "TIME" - that marker has been replaced by the boolean column above,
#    so clear it.
$ws.Cells.Item(3, 4).Value = ""

# 6) Make sure the table now spans the full A1:H5 range.
$lo.Resize($ws.Range("A1:H5"))

# 7) Cosmetic follow-up to match the edited workbook: the new column
#    keeps the (narrower) width that "variable-type" already had, and the
#    active selection on the sheet ends up on D3.
$ws.Columns.Item($newTimevalCol).ColumnWidth = $ws.Columns.Item(4).ColumnWidth
$ws.Range("D3").Select()
